$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-08 22:16:38", "Bank", "def", 10, 0.9454018282780285),
    @("2023-10-08 22:24:01", "Bank", "def", 10, 0.9454018282780285),
    @("2023-10-08 22:24:58", "Bank", "def", 10, 0.9454018282780285),
    @("2023-10-08 22:24:58", "Bank", "def", 6, 0.9456571933548517),
    @("2023-10-08 22:24:59", "Bank", "def", 7, 0.9455623965112731),
    @("2023-10-08 22:24:59", "Bank", "def", 8, 0.9454888694856882),
    @("2023-10-08 22:24:59", "Bank", "def", 5, 0.9458168575092269),
    @("2023-10-08 22:24:59", "Bank", "def", 2, 0.9518562614462791),
    @("2023-10-08 22:24:59", "Bank", "def", 3, 0.9478176730952024),
    @("2023-10-08 22:24:59", "Bank", "def", 1, 0.9600122695046023)
)

$startRow = 27
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = [string]$row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
